# Database_do_an.xlsx – schema documentation sheet update
#
# Summary of the change being reproduced (per the OOXML diff):
#   - The last mini-table ("lich_su_tinh_trang_don" / "hoa_don_chi_tiet",
#     previously at rows 25-29) gets pushed down two rows, to rows 27-31,
#     to make room for more rows in the table above it.
#   - "nhan_vien" table gains a new field: token (D11).
#   - "san_pham" table (G14 block) is re-worked: gains so_luong and
#     ma_gioi_tinh fields, keeps ma_the_loai/mo_ta/anh/ma_nsx but in a new
#     order/extent (rows 18-23 instead of 18-21).
#   - "the_loai" table (J column) drops its mo_ta / anh_dai_dien rows, so
#     it now only has ma/ten fields.
#   - A brand new small lookup table "gioi_tinh" (ma/ten) is added in
#     column L.
#   - A stray formatted (bold) but empty marker cell appears at I27,
#     alongside the header row of the last table.
#   - View state: the window is scrolled so row 10 is at the top and the
#     active selection moves to G24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the bottom table down by 2 rows -------------------------------
# Inserting two whole rows above row 25 carries rows 25-29 (with their
# values *and* bold header formatting) down to rows 27-31 automatically.
$ws.Rows("25:26").Insert()

# --- nhan_vien: add "token" field -----------------------------------------
$ws.Range("D11").Value = "token"

# --- san_pham table: re-sequence / extend its field list -------------------
$ws.Range("G18").Value = "so_luong"
$ws.Range("G19").Value = "ma_gioi_tinh"
$ws.Range("G20").Value = "ma_the_loai"
$ws.Range("G21").Value = "mo_ta"
$ws.Range("G22").Value = "anh"
$ws.Range("G23").Value = "ma_nsx"

# --- the_loai table: drop mo_ta / anh_dai_dien rows -------------------------
$ws.Range("J4:J5").ClearContents()

# --- New gioi_tinh lookup table (column L) ---------------------------------
$ws.Range("L1").Value = "gioi_tinh"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L2").Value = "ma"
$ws.Range("L3").Value = "ten"

# --- Stray bold marker cell next to the last table's header row ------------
$ws.Range("I27").Font.Bold = $true

# --- Final view/selection state --------------------------------------------
$ws.Range("G24").Select() | Out-Null
